$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values (row numbers as in the original layout) ---

# RM 14 (row 5): clear E5
$ws.Range("E5").ClearContents()

# RM 58 (row 11): set E11
$ws.Range("E11").Value = -7.9

# RM 125 (row 19): set C19, clear E19
$ws.Range("C19").Value = 13.2
$ws.Range("E19").ClearContents()

# RM 135 (row 21): clear C21
$ws.Range("C21").ClearContents()

# RM 140 (row 23): set C23, set E23
$ws.Range("C23").Value = 12.2
$ws.Range("E23").Value = -7

# RM 145 (row 25): set E25
$ws.Range("E25").Value = -7.1

# SC 5 (row 27): clear B27
$ws.Range("B27").ClearContents()

# SC 101 (row 29): set B29, clear C29, clear E29
$ws.Range("B29").Value = -20.4
$ws.Range("C29").ClearContents()
$ws.Range("E29").ClearContents()

# SC 119 (row 31): clear B31, clear E31
$ws.Range("B31").ClearContents()
$ws.Range("E31").ClearContents()

# SC 232 (row 35): set C35, set E35
$ws.Range("C35").Value = 10.4
$ws.Range("E35").Value = -10.7

# --- Remove rows entirely (delete bottom-up so earlier row numbers stay valid) ---

# SC 92 (row 28)
$ws.Rows.Item(28).Delete()

# RM 232 (row 26)
$ws.Rows.Item(26).Delete()
